# "Remove all information from control panel"
#
# The control panel sheet (A1:H3) held a labelled row of age-breakpoint
# values (A1:H1, A1 highlighted/bordered) plus two blank-but-styled helper
# rows (C2:F2 and C3:F3). The edit wipes every value while leaving the
# plain "grid" look of C1:F1 matching the already-plain C2:F2 row, and it
# drops the now-empty third row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 ("age_breakpoints" label), B1 (5) and the trailing G1:H1 cells carried
# the special header formatting (fill/border/font) - drop content and
# formatting together so they disappear from the sheet entirely.
$ws.Range("A1:B1").Clear()
$ws.Range("G1:H1").Clear()

# C1:F1 held values (15, ...) under the same highlighted header style as
# A1:B1. They should keep existing as blank cells but drop back to the
# plain look already used one row down (C2:F2), so copy that formatting
# over before clearing the leftover values.
$ws.Range("C2:F2").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1:F1").ClearContents()

# The third row (C3:F3) duplicated the plain helper row and is removed
# completely, shifting the sheet's used range up to just two rows.
$ws.Range("A3:I3").Delete()

# Leave the cursor parked just below the (now two-row) table, as in the
# edited workbook.
$ws.Range("A4").Select()
